# New crime data collected - weekly CompStat report update (121st Precinct)
# Report period moves forward one week: 12/9/2024-12/15/2024 -> 12/16/2024-12/22/2024
# Volume/Number label increments 50 -> 51, and the Crime Complaints table is
# refreshed with the new week's figures (a new "Prepared by" footer row is
# also inserted, pushing the trailing two rows down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# ---------------------------------------------------------------------------
# Crime Complaints table (rows 14-31) - refreshed counts / percentages
# ---------------------------------------------------------------------------

# Row 14 - Murder : 28-day count becomes unavailable (N/A)
$ws.Range("F14").Value = "'0"

# Row 15 - Rape
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 50

# Row 16 - Robbery
$ws.Range("D16").Value = 2
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -71.428571428571
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 7.5
$ws.Range("L16").Value = 3.614457831325

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 35.714285714285
$ws.Range("I17").Value = 203
$ws.Range("J17").Value = 242
$ws.Range("K17").Value = -16.115702479338
$ws.Range("L17").Value = -4.694835680751

# Row 18 - Burglary
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 38.571428571428
$ws.Range("L18").Value = 8.988764044943

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 350
$ws.Range("J19").Value = 414
$ws.Range("K19").Value = -15.458937198067
$ws.Range("L19").Value = -15.048543689320

# Row 20 - G.L.A.
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -60
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = -18.348623853211
$ws.Range("L20").Value = -19.819819819819

# Row 21 - TOTAL
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -23.809523809523
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = -7.575757575757
$ws.Range("I21").Value = 847
$ws.Range("J21").Value = 931
$ws.Range("K21").Value = -9.022556390977
$ws.Range("L21").Value = -8.034744842562

# Row 23 - Housing
$ws.Range("L23").Value = -15

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 51.612903225806
$ws.Range("F24").Value = 156
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 1458
$ws.Range("J24").Value = 1459
$ws.Range("K24").Value = -0.068540095956
$ws.Range("L24").Value = 6.190823015294

# Row 25 - Retail Theft
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 94.117647058823
$ws.Range("F25").Value = 105
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = 94.444444444444
$ws.Range("I25").Value = 992
$ws.Range("J25").Value = 839
$ws.Range("K25").Value = 18.235995232419
$ws.Range("L25").Value = 43.768115942029

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -13.157894736842
$ws.Range("I26").Value = 562
$ws.Range("J26").Value = 504
$ws.Range("K26").Value = 11.507936507936
$ws.Range("L26").Value = 16.115702479338

# Row 27 - UCR Rape*
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 20.833333333333
$ws.Range("L27").Value = 38.095238095238

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("I28").Value = 59
$ws.Range("K28").Value = -4.838709677419
$ws.Range("L28").Value = 40.476190476190

# Row 31 - Hate Crimes : this week now unavailable (N/A)
$ws.Range("C31").Value = "'0"
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "***.*"

# ---------------------------------------------------------------------------
# Insert the new "Prepared by" footer row (pushes old rows 56-57 to 57-58)
# ---------------------------------------------------------------------------
$ws.Rows(56).Insert()
